$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.255.06"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.547.86"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.42"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.53"
$ws.Range("E6").Value = "  +5.40%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.547.19"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("E12").Value = "  -3.02%  "

$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.04"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.012.40"
$ws.Range("E15").Value = "  -2.38%  "

$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.107.25"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.551.14"
$ws.Range("E18").Value = "  -2.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.05"
$ws.Range("E19").Value = "  +2.84%  "

$ws.Range("E20").Value = "  -3.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "355.53"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +3.90%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.04"
$ws.Range("E26").Value = "  +1.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -4.81%  "

$ws.Range("E29").Value = "  +0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0996"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "535.64"
$ws.Range("E31").Value = "  -1.54%  "

$ws.Range("E32").Value = "  +4.70%  "

$ws.Range("E33").Value = "  +0.64%  "

$ws.Range("E34").Value = "  -0.64%  "

$ws.Range("E35").Value = "  -1.55%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.44"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.81"
$ws.Range("E39").Value = "  -0.66%  "

$ws.Range("E40").Value = "  +1.09%  "

$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  +1.11%  "

$ws.Range("E44").Value = "  +6.82%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.82"
$ws.Range("E46").Value = "  -0.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.93"
$ws.Range("E47").Value = "  +0.23%  "

$ws.Range("E48").Value = "  -1.81%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0281"
$ws.Range("E49").Value = "  -6.11%  "

$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("E51").Value = "  +1.21%  "

